$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.335.19'
$ws.Range('E2').Value = '  +1.11%  '
$ws.Range('D3').Value = '1.856.75'
$ws.Range('E3').Value = '  +1.53%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.43%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '314.03'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.02%  '
$ws.Range('E6').Value = '  -0.53%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4617'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.15%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3698'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.26%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07319'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.83%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8822'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.03%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07828'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.21%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '19.81'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.86%  '
$ws.Range('D13').Value = '1.846.50'
$ws.Range('E13').Value = '  -0.12%  '
$ws.Range('E14').Value = '  +0.59%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.535'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.06%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '91.81'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.17%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.003'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.43%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008858'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.76%  '
$ws.Range('E19').Value = '  -0.61%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.81'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.68%  '
$ws.Range('D21').Value = '27.357.57'
$ws.Range('E21').Value = '  +0.77%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.112'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.00%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.49'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.73%  '
$ws.Range('D24').Value = '2.080.34'
$ws.Range('E24').Value = '  +0.15%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.886'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.52%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '152.03'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.87%  '
$ws.Range('E27').Value = '  +0.83%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.071'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.23%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.115'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.38%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '115.92'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.42%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08848'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.00%  '
$ws.Range('E32').Value = '  +4.86%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.999'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.02%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.168'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.77%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.487'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.96%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.606'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +4.56%  '
$ws.Range('E37').Value = '  +0.90%  '
$ws.Range('E38').Value = '  +0.11%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.990'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.05%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.05215'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.48%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '7.025'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.99%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.5146'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.55%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1637'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.48%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.338'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.54%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.4831'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.12%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.30'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.00%  '
$ws.Range('E47').Value = '  -0.64%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '103.06'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.22%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.650'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.74%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06222'
$ws.Range('D50').Style = 'Normal'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '65.55'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.87%  '
